# Weekly update: insert a new observation row for
# "Femacal de La Calera - Zanahoria" dated 2021-08-30 (serial 44438),
# pushing all existing rows (113 onward) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 113; this shifts rows 113:181 down to 114:182
# (values, formats and styles move with them automatically).
$ws.Rows.Item(113).Insert()

# Populate the newly inserted row 113 with the new weekly record.
$ws.Cells.Item(113, 1).Value  = 3
$ws.Cells.Item(113, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(113, 3).Value  = "Coquimbo"
$ws.Cells.Item(113, 4).Value  = 44438
$ws.Cells.Item(113, 5).Value  = 5
$ws.Cells.Item(113, 6).Value  = 100114013
$ws.Cells.Item(113, 7).Value  = "Zanahoria"
$ws.Cells.Item(113, 8).Value  = "Sin especificar"
$ws.Cells.Item(113, 9).Value  = "Primera"
$ws.Cells.Item(113, 10).Value = 440
$ws.Cells.Item(113, 11).Value = 5000
$ws.Cells.Item(113, 12).Value = 5500
$ws.Cells.Item(113, 13).Value = 5205
$ws.Cells.Item(113, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(113, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(113, 16).Value = 260
$ws.Cells.Item(113, 17).Value = 20
$ws.Cells.Item(113, 18).Value = "Hortaliza"
